# Update "want-to-go" counts (column F) on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 194
$ws1.Range("F8").Value  = 10079
$ws1.Range("F10").Value = 3468
$ws1.Range("F14").Value = 2767
$ws1.Range("F23").Value = 129
$ws1.Range("F27").Value = 613
$ws1.Range("F34").Value = 2754
$ws1.Range("F35").Value = 2950
$ws1.Range("F42").Value = 86

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 194
$ws4.Range("F11").Value = 10079
$ws4.Range("F18").Value = 2767
$ws4.Range("F24").Value = 129
$ws4.Range("F27").Value = 613
$ws4.Range("F33").Value = 2754
$ws4.Range("F35").Value = 2950
$ws4.Range("F45").Value = 86
